$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.213.90'
$ws.Range('E2').Value = '  -2.15%  '
$ws.Range('D3').Value = '3.754.96'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.68'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.80'
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('D7').Value = '3.752.67'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('E10').Value = '  -3.07%  '
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  -4.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.74'
$ws.Range('E14').Value = '  -3.22%  '
$ws.Range('D15').Value = '4.385.85'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '3.774.93'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '67.165.19'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.68'
$ws.Range('E18').Value = '  -2.85%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.92'
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.46'
$ws.Range('E21').Value = '  -5.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '456.90'
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.694'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000150'
$ws.Range('E24').Value = '  +4.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.87'
$ws.Range('E25').Value = '  -2.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.78'
$ws.Range('E26').Value = '  -3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.11'
$ws.Range('E27').Value = '  -5.85%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.96'
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.76'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.19'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '29.58'
$ws.Range('E32').Value = '  -2.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  -3.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.18'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '3.709.30'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0993'
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  -6.07%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.137'
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.990'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.72'
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.71'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.297'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.68'
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.34'
$ws.Range('E47').Value = '  -3.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '146.80'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.82'
$ws.Range('E49').Value = '  -7.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '387.25'
$ws.Range('E50').Value = '  -3.40%  '
$ws.Range('D51').Value = '2.742.26'
$ws.Range('E51').Value = '  +1.75%  '
